$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: extend the separator row's bottom-border style into the new column
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Row 3: new year header "2021" - same look as the other year headers (H3)
# but bumped to 11pt (matches the style/font additions in the diff)
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 2021
$ws.Range("I3").Font.Size = 11

# Row 4: new data value 149 - same look as H4, bumped to 11pt
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 149
$ws.Range("I4").Font.Size = 11

# Row 5: new data value 159 - same look as H5, bumped to 11pt
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 159
$ws.Range("I5").Font.Size = 11

# Match the saved selection state from the diff
$ws.Range("K4").Select() | Out-Null
